$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 with a single cell A10 = "." using the same formatting as
# the other "Alunos" column cells (e.g. A2, A5, A6, A7, A8).
$ws.Range("A10").Value = "."

# Match font formatting of the reference cell (Arial 10, black).
$ws.Range("A10").Font.Name = $ws.Range("A2").Font.Name
$ws.Range("A10").Font.Size = $ws.Range("A2").Font.Size
$ws.Range("A10").Font.Color = $ws.Range("A2").Font.Color
$ws.Range("A10").Font.Bold = $ws.Range("A2").Font.Bold
$ws.Range("A10").Font.Italic = $ws.Range("A2").Font.Italic
$ws.Range("A10").Font.Underline = $ws.Range("A2").Font.Underline

# Match the row height/format used by the other data rows.
$ws.Rows.Item(10).RowHeight = $ws.Rows.Item(2).RowHeight

# Update the active selection to A10 (matches the new selection in the diff).
$ws.Range("A10").Select()
